# The sheet holds 4 years (2014-2017) of monthly data in rows 2..49, one
# row per month, in chronological Jan->Dec order for each year. The edit
# re-sorts each year's 12-month block so that Oct, Nov, Dec come first,
# followed by Jan..Sep (i.e. a "move last quarter to the top" re-sort),
# while leaving the header row (row 1) and the underlying values
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 49

# Read the existing data rows into parallel arrays.
$months = @()
$valB = @()
$valC = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $months += $ws.Cells.Item($r, 1).Value2
    $valB += $ws.Cells.Item($r, 2).Value2
    $valC += $ws.Cells.Item($r, 3).Value2
}

# Re-order: within every consecutive block of 12 rows (one calendar year),
# move the last 3 entries (Oct, Nov, Dec) to the front of that block, ahead
# of Jan..Sep.
$newMonths = @()
$newB = @()
$newC = @()
$rowCount = $lastRow - $firstRow + 1
$yearCount = $rowCount / 12
for ($y = 0; $y -lt $yearCount; $y++) {
    $base = $y * 12

    for ($k = 9; $k -le 11; $k++) {
        $newMonths += $months[$base + $k]
        $newB += $valB[$base + $k]
        $newC += $valC[$base + $k]
    }

    for ($k = 0; $k -le 8; $k++) {
        $newMonths += $months[$base + $k]
        $newB += $valB[$base + $k]
        $newC += $valC[$base + $k]
    }
}

# Write the re-ordered data back in place.
for ($i = 0; $i -lt $newMonths.Count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value = $newMonths[$i]
    $ws.Cells.Item($r, 2).Value = $newB[$i]
    $ws.Cells.Item($r, 3).Value = $newC[$i]
}
